# Apply cell-value edits to the "Лист1" worksheet (gradebook-style sheet,
# columns B..I are raw scores, J/K/L are derived formulas that recalc
# automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# Row 3 (Архипова Светлана Владимировна): D3 0 -> 5
$ws.Range("D3").Value = 5

# Row 4 (Быков Вадим Дмитриевич): E4 empty -> 5
$ws.Range("E4").Value = 5

# Row 5 (Галямова Яна Дмитриевна): B5 empty -> -2, D5 empty -> -1
$ws.Range("B5").Value = -2
$ws.Range("D5").Value = -1

# Row 7 (Жамсаранова Аяна Жаргаловна): E7 empty -> 5
$ws.Range("E7").Value = 5

# Row 8 (Захаренкова Екатерина Денисовна): B8 empty -> -2, D8 empty -> -1
$ws.Range("B8").Value = -2
$ws.Range("D8").Value = -1

# Row 10 (Косарынская Анна Александровна): E10 empty -> 5
$ws.Range("E10").Value = 5

# Row 12 (Круглов Кирилл Максимович): D12 empty -> 4
$ws.Range("D12").Value = 4

# Row 13 (Крутов Никита Сергеевич): D13 0 -> 5
$ws.Range("D13").Value = 5

# Row 14 (Ларюшин Виктор Романович): D14 0 -> 5
$ws.Range("D14").Value = 5

# Row 16 (Мачкалян Тигран Норайрович): G16 empty -> 5, H16 empty -> 5
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 5

# Row 17 (Молокова Татьяна Михайловна): C17 2 -> 3
$ws.Range("C17").Value = 3

# Row 20 (Ротанкова Вера Владимировна): E20 empty -> 5
$ws.Range("E20").Value = 5

# Row 22 (Сычиков Владимир Андреевич): B22 empty -> 3, D22 empty -> -1
$ws.Range("B22").Value = 3
$ws.Range("D22").Value = -1

# Row 23 (Ушакова Александра Юрьевна): E23 empty -> 5
$ws.Range("E23").Value = 5

# Row 26 (Шаблыгин Михаил Максимович): E26 0 -> 5
$ws.Range("E26").Value = 5

# Update the saved selection / active cell to B9, matching the author's
# last cursor position when the file was uploaded.
$ws.Activate()
$ws.Range("B9").Select()
